# Update the cryptocurrency price (col D) and 1h volume change (col E)
# snapshot for the "cryptos" sheet, per the latest GitHub Actions refresh.
#
# Note: several "Price" values are plain numeric-looking strings (e.g.
# "97.22"). Assigning those directly to .Value would make Excel coerce
# them to real numbers (losing the original fixed-precision text/trailing
# zeros). Prefixing with a leading apostrophe forces Excel to keep them
# as text, matching the original cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.682.82"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.556.58"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'302.34"
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("D6").Value = "'97.22"
$ws.Range("E6").Value = "  +7.01%  "
$ws.Range("D7").Value = "'0.573"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "'36.48"
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").Value = "  +9.51%  "
$ws.Range("D13").Value = "'7.51"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").Value = "2.541.64"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "'14.48"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("D17").Value = "42.727.62"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "'13.62"
$ws.Range("E18").Value = "  +9.25%  "
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("D20").Value = "'6.58"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").Value = "'71.57"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "'256.04"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").Value = "'2.96"
$ws.Range("E23").Value = "  +2.62%  "
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").Value = "'28.04"
$ws.Range("E25").Value = "  -4.11%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'39.09"
$ws.Range("E27").Value = "  +8.85%  "
$ws.Range("D28").Value = "'10.02"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("D31").Value = "'155.74"
$ws.Range("E31").Value = "  +3.54%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").Value = "'0.0802"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("D36").Value = "'26.26"
$ws.Range("E36").Value = "  +9.35%  "
$ws.Range("D37").Value = "'18.33"
$ws.Range("E37").Value = "  +16.85%  "
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("E40").Value = "  +1.56%  "
$ws.Range("E41").Value = "  +29.88%  "
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "2.061.24"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'88.01"
$ws.Range("E46").Value = "  +4.67%  "
$ws.Range("D47").Value = "'9.21"
$ws.Range("E47").Value = "  +5.79%  "
$ws.Range("D48").Value = "'76.34"
$ws.Range("E48").Value = "  +10.92%  "
$ws.Range("D49").Value = "2.805.72"
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").Value = "'103.67"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "'0.189"
$ws.Range("E51").Value = "  +2.88%  "
